$d = $word.ActiveDocument

function Find-ParagraphByPrefix($doc, $prefix) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.StartsWith($prefix)) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) Remove the "Hamã, Hebreu, ..." key-terms listing paragraph (it sits right
#    after the "H" Heading2 paragraph).
# ---------------------------------------------------------------------------
$pHama = Find-ParagraphByPrefix $d "Hamã, Hebreu"
$pHama.Range.Delete()

# ---------------------------------------------------------------------------
# 2) Remove the "This PDF version is provided under the same license."
#    paragraph entirely.
# ---------------------------------------------------------------------------
$pPdf = Find-ParagraphByPrefix $d "This PDF version is provided under the same license."
$pPdf.Range.Delete()

# ---------------------------------------------------------------------------
# 3) Rewrite the license-description paragraph:
#    "Termos Chave (Biblica) (Portuguese) is based on: Biblica Bible
#    Dictionary, Biblica, Inc., 2023, which is licensed under a CC BY-SA 4.0
#    license."
#    becomes a new description of "Biblica Study Notes (Key Terms)", dropping
#    the two hyperlinks and adding the adapted-languages sentence.
# ---------------------------------------------------------------------------
$pLicense = Find-ParagraphByPrefix $d "Termos Chave (Biblica) (Portuguese) is based on"

# 3a) Rename the bold title run. This text sits at the very start of the
#     paragraph so there is nothing to its left for Word to merge into.
$scope = $d.Range($pLicense.Range.Start, $pLicense.Range.End)
$scope.Find.Execute("Termos Chave (Biblica)", $true, $false, $false, $false, $false, $true, 1, $false, "Biblica Study Notes (Key Terms)", 2) | Out-Null

# 3b) Replace the remainder of the paragraph (everything after the bold run,
#     up to the trailing empty run/paragraph mark) with the new license +
#     adaptation text. The leading space is left out of the search match so
#     the replacement does not get absorbed into the preceding bold run.
$pLicense = Find-ParagraphByPrefix $d "Biblica Study Notes (Key Terms)"
$scope = $d.Range($pLicense.Range.Start, $pLicense.Range.End)
$oldTail = "(Portuguese) is based on: Biblica Bible Dictionary, Biblica, Inc., 2023, which is licensed under a CC BY-SA 4.0 license."
$newTail = "© 2023 Biblica Inc. Released under CC BY-SA 4.0 license. Biblica Study Notes has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文)from Biblica Study Notes © 2023 Biblica Inc. Released under CC BY-SA 4.0 license by Mission Mutual."
$scope.Find.Execute($oldTail, $true, $false, $false, $false, $false, $true, 1, $false, $newTail, 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Remove the "License Information" Heading2 paragraph entirely.
# ---------------------------------------------------------------------------
$pInfo = Find-ParagraphByPrefix $d "License Information"
$pInfo.Range.Delete()
